$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 306.9
$ws.Range("I6").Value = 159.33333
$ws.Range("K6").Value = 477.99999
$ws.Range("M6").Value = -365.99999
$ws.Range("H8").Value = 267.71155
$ws.Range("I8").Value = 51.75
$ws.Range("K8").Value = 155.25
$ws.Range("M8").Value = -16.25
$ws.Range("H9").Value = 646.8570999999999
$ws.Range("I9").Value = 810
$ws.Range("J9").Value = 524.5
$ws.Range("K9").Value = 810
$ws.Range("L9").Value = 524.5
$ws.Range("M9").Value = -641
$ws.Range("N9").Value = -862.5
$ws.Range("H58").Value = 45454932
$ws.Range("J58").Value = 166667630
$ws.Range("L58").Value = 500002890
$ws.Range("N58").Value = -500003190
$ws.Range("H103").Value = 1095.3334
$ws.Range("I103").Value = 900.375
$ws.Range("J103").Value = 1251.3
$ws.Range("K103").Value = 2701.125
$ws.Range("L103").Value = 3753.9
$ws.Range("M103").Value = -2115.125
$ws.Range("N103").Value = -4925.9
$ws.Range("H123").Value = 54995.6
$ws.Range("J123").Value = 54995.6
$ws.Range("L123").Value = 54995.6
$ws.Range("N123").Value = -64795.6
$ws.Range("H132").Value = 1052
$ws.Range("I132").Value = 1052
$ws.Range("K132").Value = 3156
$ws.Range("M132").Value = -626
$ws.Range("H135").Value = 2000505.9
$ws.Range("I135").Value = 2500378
$ws.Range("K135").Value = 22503402
$ws.Range("M135").Value = -22500867

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3853038
$ws.Range("I32").Value = 4171853.8
$ws.Range("K32").Value = 4171853.8
$ws.Range("M32").Value = -4171566.8
$ws.Range("H34").Value = 87116.92
$ws.Range("I34").Value = 180008.33
$ws.Range("K34").Value = 180008.33
$ws.Range("M34").Value = -179737.33
$ws.Range("H74").Value = 33780.906
$ws.Range("I74").Value = 41599
$ws.Range("K74").Value = 41599
$ws.Range("M74").Value = -40725
$ws.Range("H77").Value = 33780.906
$ws.Range("I77").Value = 41599
$ws.Range("K77").Value = 207995
$ws.Range("M77").Value = -203627
$ws.Range("H97").Value = 2452308.8
$ws.Range("I97").Value = 1204.24
$ws.Range("J97").Value = 9260932
$ws.Range("K97").Value = 1204.24
$ws.Range("L97").Value = 9260932
$ws.Range("M97").Value = -708.24
$ws.Range("N97").Value = -9261924
$ws.Range("H122").Value = 4534.2334
$ws.Range("I122").Value = 3193.2273
$ws.Range("J122").Value = 8222
$ws.Range("K122").Value = 9579.6819
$ws.Range("L122").Value = 24666
$ws.Range("M122").Value = -7129.6819
$ws.Range("N122").Value = -29566
$ws.Range("H132").Value = 3810.8115
$ws.Range("I132").Value = 2954.647
$ws.Range("K132").Value = 8863.940999999999
$ws.Range("M132").Value = -6333.940999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1811.7435
$ws.Range("I94").Value = 756.8
$ws.Range("J94").Value = 5328.222
$ws.Range("K94").Value = 756.8
$ws.Range("L94").Value = 5328.222
$ws.Range("M94").Value = -305.8
$ws.Range("N94").Value = -6230.222
$ws.Range("H134").Value = 5322483.5
$ws.Range("I134").Value = 8065970.5
$ws.Range("J134").Value = 6977.75
$ws.Range("K134").Value = 24197911.5
$ws.Range("L134").Value = 20933.25
$ws.Range("M134").Value = -24195376.5
$ws.Range("N134").Value = -26003.25
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 74890
$ws.Range("J137").Value = 74890
$ws.Range("L137").Value = 74890
$ws.Range("N137").Value = -85090

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 69193.25
$ws.Range("J59").Value = 69193.25
$ws.Range("L59").Value = 69193.25
$ws.Range("N59").Value = -71483.25
$ws.Range("H62").Value = 4824.7393
$ws.Range("J62").Value = 5413.5713
$ws.Range("L62").Value = 5413.5713
$ws.Range("N62").Value = -6661.5713
$ws.Range("H65").Value = 4824.7393
$ws.Range("J65").Value = 5413.5713
$ws.Range("L65").Value = 27067.8565
$ws.Range("N65").Value = -33307.85649999999
$ws.Range("H122").Value = 4836.587
$ws.Range("I122").Value = 4225.657
$ws.Range("J122").Value = 6780.4546
$ws.Range("K122").Value = 12676.971
$ws.Range("L122").Value = 20341.3638
$ws.Range("M122").Value = -10226.971
$ws.Range("N122").Value = -25241.3638
$ws.Range("H132").Value = 3392.9111
$ws.Range("I132").Value = 2582.125
$ws.Range("J132").Value = 5388.6924
$ws.Range("K132").Value = 7746.375
$ws.Range("L132").Value = 16166.0772
$ws.Range("M132").Value = -5216.375
$ws.Range("N132").Value = -21226.0772
$ws.Range("H134").Value = 3494.5435
$ws.Range("I134").Value = 2210.5312
$ws.Range("K134").Value = 6631.5936
$ws.Range("M134").Value = -4096.5936

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 22731452
$ws.Range("I56").Value = 22731452
$ws.Range("K56").Value = 22731452
$ws.Range("M56").Value = -22730922
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 9000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -8314
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 27000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -23568
$ws.Range("N65").ClearContents()
$ws.Range("H87").Value = 606
$ws.Range("I87").Value = 606
$ws.Range("K87").Value = 1818
$ws.Range("M87").Value = -570
$ws.Range("H90").Value = 606
$ws.Range("I90").Value = 606
$ws.Range("K90").Value = 5454
$ws.Range("M90").Value = 786
$ws.Range("H107").Value = 15385422
$ws.Range("J107").Value = 20000894
$ws.Range("L107").Value = 60002682
$ws.Range("N107").Value = -60006522

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 83749.75
$ws.Range("J52").Value = 90000
$ws.Range("L52").Value = 90000
$ws.Range("N52").Value = -90518
$ws.Range("H70").Value = 8424.214
$ws.Range("I70").Value = 4826.5
$ws.Range("K70").Value = 4826.5
$ws.Range("M70").Value = -4556.5
$ws.Range("H73").Value = 8424.214
$ws.Range("I73").Value = 4826.5
$ws.Range("K73").Value = 4826.5
$ws.Range("M73").Value = -3890.5
$ws.Range("H113").Value = 5431.952
$ws.Range("I113").Value = 2158.125
$ws.Range("J113").Value = 9797.056
$ws.Range("K113").Value = 2158.125
$ws.Range("L113").Value = 9797.056
$ws.Range("M113").Value = 11.875
$ws.Range("N113").Value = -14137.056
$ws.Range("H132").Value = 1672.025
$ws.Range("I132").Value = 1222.5588
$ws.Range("K132").Value = 3667.6764
$ws.Range("M132").Value = -1137.6764

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 12513749
$ws.Range("I33").Value = 25008748
$ws.Range("J33").Value = 18750
$ws.Range("K33").Value = 25008748
$ws.Range("L33").Value = 18750
$ws.Range("M33").Value = -25008458
$ws.Range("N33").Value = -19330
$ws.Range("H40").Value = 20836686
$ws.Range("I40").Value = 26317920
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 26317920
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -26317784
$ws.Range("N40").Value = -8272
$ws.Range("H69").Value = 53963
$ws.Range("J69").Value = 53963
$ws.Range("L69").Value = 53963
$ws.Range("N69").Value = -55585
$ws.Range("H72").Value = 53963
$ws.Range("J72").Value = 53963
$ws.Range("L72").Value = 161889
$ws.Range("N72").Value = -170001
$ws.Range("H80").Value = 49993.332
$ws.Range("J80").Value = 49993.332
$ws.Range("L80").Value = 49993.332
$ws.Range("N80").Value = -52239.332
$ws.Range("H83").Value = 49993.332
$ws.Range("J83").Value = 49993.332
$ws.Range("L83").Value = 149979.996
$ws.Range("N83").Value = -161211.996
$ws.Range("H136").Value = 8280.725
$ws.Range("I136").Value = 2981.32
$ws.Range("J136").Value = 12295.424
$ws.Range("K136").Value = 8943.960000000001
$ws.Range("L136").Value = 36886.272
$ws.Range("M136").Value = -6393.960000000001
$ws.Range("N136").Value = -41986.272

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 22914620
$ws.Range("I122").Value = 36003588
$ws.Range("K122").Value = 108010764
$ws.Range("M122").Value = -108008314
$ws.Range("H132").Value = 4950.8613
$ws.Range("I132").Value = 4730.6294
$ws.Range("K132").Value = 14191.8882
$ws.Range("M132").Value = -11661.8882
